$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- N15: add "Sprites" text ---
$ws.Range("N15").Value = "Sprites"

# --- Row 16 second table: fill in J16:L16, N16 ---
$ws.Range("J16").Value = 45391
$ws.Range("K16").Formula = "=8"
$ws.Range("L16").Formula = "=9+35/60"
$ws.Range("N16").Value = "Sprites"

# --- New row 34 ---
$ws.Range("A34").Value = 45391
$ws.Range("B34").Formula = "=8"
$ws.Range("C34").Formula = "=9+35/60"
$ws.Range("D34").Formula = "=C34-B34"
$ws.Range("E34").Value = "Level designs implementieren"
